$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.02"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.74"
$ws.Range("D3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.199"
$ws.Range("D4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06096"
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.514"
$ws.Range("D6").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1578"
$ws.Range("D10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08090"
$ws.Range("D11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03344"
$ws.Range("D12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03121"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09283"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.919"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001689"
$ws.Range("D16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04814"
$ws.Range("D17").ClearFormats()

$ws.Range("B18").Value = "TigerCash"

$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006192"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "BitKan"

$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.001101"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"

$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.003392"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "19HotbitTokenHTBWorstin24h"

$ws.Range("B21").Value = "NitroEx"

$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001499"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"

$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.695"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"

$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.263"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"

$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01328"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "23OneONE"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3358"
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1275"
$ws.Range("D26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006161"
$ws.Range("D27").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04609"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007185"
$ws.Range("D41").ClearFormats()

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003898"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1121"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.002969"
$ws.Range("D45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006026"
$ws.Range("D46").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7495"
$ws.Range("D48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1250"
$ws.Range("D49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001499"
$ws.Range("D50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01009"
$ws.Range("D51").ClearFormats()
